# Scheduler input workbook update
#
# Workbook layout:
#  - Worksheet "Sheet1" is the course-data table (columns A:H, rows 1-12) and
#    is the selected/active tab.
#  - Worksheet "Sheet2" is the "Conflict Courses" reference table
#    (columns A:C, rows 1-12).

$wb = $excel.ActiveWorkbook

$data      = $wb.Worksheets.Item("Sheet1")
$conflicts = $wb.Worksheets.Item("Sheet2")

# 1) The "num of sections" column (D) becomes 1 for every course row (2-12)
for ($r = 2; $r -le 12; $r++) {
    $data.Cells.Item($r, 4).Value = 1
}

# 2) Every course now shares the same "instructor hours" session window, so
#    retype column H (rows 4-12) to the new time of 08:00 / 17:30.
for ($r = 4; $r -le 12; $r++) {
    $data.Cells.Item($r, 8).Value = "08:00 / 17:30"
}

# 3) LANG202 (row 12) now only meets on Friday instead of Monday / Wednesday
$data.Range("G12").Value = "Friday"

# 4) Leave the same cell selections behind on each sheet
$data.Range("G13").Select()
$conflicts.Range("C3").Select()

# Keep the data sheet ("Sheet1") as the active/selected tab
$data.Activate()
